# ============================================================================
# Scheduled-runner refresh: pulls current Market Board prices for each Leve
# crafting recipe and rewrites the derived profit columns (H:N) on every job
# sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Row/column layout is
# unchanged; only the numeric price + profit figures are refreshed here.
#   H = currentAveragePrice        K = LevePriceNQ
#   I = currentAveragePriceNQ      L = LevePriceHQ
#   J = currentAveragePriceHQ      M = LeveProfitNQ   N = LeveProfitHQ
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 149750
$ws.Range("I18").Value = 1000
$ws.Range("K18").Value = 1000
$ws.Range("M18").Value = -716
$ws.Range("H40").Value = 4087.25
$ws.Range("I40").Value = 3399.6667
$ws.Range("K40").Value = 3399.6667
$ws.Range("M40").Value = -3224.6667
$ws.Range("H43").Value = 6610.4287
$ws.Range("I43").Value = 6880
$ws.Range("J43").Value = 6460.6665
$ws.Range("K43").Value = 6880
$ws.Range("L43").Value = 6460.6665
$ws.Range("M43").Value = -6811
$ws.Range("N43").Value = -6598.6665
$ws.Range("H96").Value = 815.375
$ws.Range("I96").Value = 665
$ws.Range("K96").Value = 1995
$ws.Range("M96").Value = -622
$ws.Range("H132").Value = 2563.2593
$ws.Range("I132").Value = 1405.3478
$ws.Range("K132").Value = 4216.0434
$ws.Range("M132").Value = -1686.0434
$ws.Range("H135").Value = 802.9524
$ws.Range("I135").Value = 544.5294
$ws.Range("K135").Value = 4900.7646
$ws.Range("M135").Value = -2365.7646
$ws.Range("H138").Value = 2307.8333
$ws.Range("I138").Value = 1374.1
$ws.Range("J138").Value = 4175.3
$ws.Range("K138").Value = 4122.299999999999
$ws.Range("L138").Value = 12525.9
$ws.Range("M138").Value = 1017.700000000001
$ws.Range("N138").Value = -22805.9

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 234.07143
$ws.Range("I5").Value = 154.14285
$ws.Range("J5").Value = 314
$ws.Range("K5").Value = 154.14285
$ws.Range("L5").Value = 314
$ws.Range("M5").Value = -42.14285000000001
$ws.Range("N5").Value = -538
$ws.Range("H32").Value = 30369.473
$ws.Range("I32").Value = 16372.924
$ws.Range("K32").Value = 16372.924
$ws.Range("M32").Value = -16085.924
$ws.Range("H63").Value = 2390.4
$ws.Range("I63").Value = 2390.4
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2390.4
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1704.4
$ws.Range("N63").Value = $null
$ws.Range("H66").Value = 2390.4
$ws.Range("I66").Value = 2390.4
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 11952
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -8520
$ws.Range("N66").Value = $null
$ws.Range("H74").Value = 1943.7778
$ws.Range("I74").Value = 1612.7333
$ws.Range("J74").Value = 3599
$ws.Range("K74").Value = 1612.7333
$ws.Range("L74").Value = 3599
$ws.Range("M74").Value = -738.7333000000001
$ws.Range("N74").Value = -5347
$ws.Range("H77").Value = 1943.7778
$ws.Range("I77").Value = 1612.7333
$ws.Range("J77").Value = 3599
$ws.Range("K77").Value = 8063.6665
$ws.Range("L77").Value = 17995
$ws.Range("M77").Value = -3695.6665
$ws.Range("N77").Value = -26731
$ws.Range("H132").Value = 2030.8846
$ws.Range("I132").Value = 1520
$ws.Range("K132").Value = 4560
$ws.Range("M132").Value = -2030

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 234.07143
$ws.Range("I4").Value = 154.14285
$ws.Range("J4").Value = 314
$ws.Range("K4").Value = 154.14285
$ws.Range("L4").Value = 314
$ws.Range("M4").Value = -39.14285000000001
$ws.Range("N4").Value = -544
$ws.Range("H96").Value = 30106.75
$ws.Range("J96").Value = 49999.5
$ws.Range("L96").Value = 49999.5
$ws.Range("N96").Value = -55491.5
$ws.Range("H134").Value = 966.4815
$ws.Range("I134").Value = 966.4815
$ws.Range("K134").Value = 2899.4445
$ws.Range("M134").Value = -364.4445000000001

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 35714652
$ws.Range("I7").Value = 50000268
$ws.Range("K7").Value = 50000268
$ws.Range("M7").Value = -50000155
$ws.Range("H58").Value = 899
$ws.Range("I58").Value = 928.2941
$ws.Range("J58").Value = 733
$ws.Range("K58").Value = 928.2941
$ws.Range("L58").Value = 733
$ws.Range("M58").Value = -725.2941
$ws.Range("N58").Value = -1139
$ws.Range("H132").Value = 1792.122
$ws.Range("I132").Value = 1618.3684
$ws.Range("K132").Value = 4855.1052
$ws.Range("M132").Value = -2325.1052
$ws.Range("H134").Value = 1715.6552
$ws.Range("I134").Value = 1885.5217
$ws.Range("J134").Value = 1064.5
$ws.Range("K134").Value = 5656.5651
$ws.Range("L134").Value = 3193.5
$ws.Range("M134").Value = -3121.5651
$ws.Range("N134").Value = -8263.5
$ws.Range("H136").Value = 899
$ws.Range("I136").Value = 928.2941
$ws.Range("J136").Value = 733
$ws.Range("K136").Value = 2784.8823
$ws.Range("L136").Value = 2199
$ws.Range("M136").Value = -234.8822999999998
$ws.Range("N136").Value = -7299
$ws.Range("H141").Value = 113941.5
$ws.Range("J141").Value = 113941.5
$ws.Range("L141").Value = 113941.5
$ws.Range("N141").Value = -124301.5

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 662
$ws.Range("I5").Value = 573.2
$ws.Range("K5").Value = 1719.6
$ws.Range("M5").Value = -1607.6
$ws.Range("H70").Value = 5066.6665
$ws.Range("I70").Value = 3725
$ws.Range("J70").Value = 7750
$ws.Range("K70").Value = 11175
$ws.Range("L70").Value = 23250
$ws.Range("M70").Value = -10860
$ws.Range("N70").Value = -23880
$ws.Range("H73").Value = 5066.6665
$ws.Range("I73").Value = 3725
$ws.Range("J73").Value = 7750
$ws.Range("K73").Value = 11175
$ws.Range("L73").Value = 23250
$ws.Range("M73").Value = -10083
$ws.Range("N73").Value = -25434
$ws.Range("H109").Value = 2087.5
$ws.Range("I109").Value = 2087.5
$ws.Range("K109").Value = 6262.5
$ws.Range("M109").Value = -5222.5
$ws.Range("H128").Value = 399660.66
$ws.Range("I128").Value = 399660.66
$ws.Range("K128").Value = 1198981.98
$ws.Range("M128").Value = -1194001.98
$ws.Range("H135").Value = 662
$ws.Range("I135").Value = 573.2
$ws.Range("K135").Value = 5158.8
$ws.Range("M135").Value = -2623.8

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4693.5293
$ws.Range("J80").Value = 5898.875
$ws.Range("L80").Value = 5898.875
$ws.Range("N80").Value = -7894.875
$ws.Range("H83").Value = 4693.5293
$ws.Range("J83").Value = 5898.875
$ws.Range("L83").Value = 29494.375
$ws.Range("N83").Value = -39478.375
$ws.Range("H93").Value = 49999.668
$ws.Range("J93").Value = 49999.668
$ws.Range("L93").Value = 49999.668
$ws.Range("N93").Value = -53743.668
$ws.Range("H122").Value = 133446.31
$ws.Range("I122").Value = 144733.34
$ws.Range("J122").Value = 1764.3334
$ws.Range("K122").Value = 434200.02
$ws.Range("L122").Value = 5293.0002
$ws.Range("M122").Value = -431750.02
$ws.Range("N122").Value = -10193.0002

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1486.9375
$ws.Range("I46").Value = 1758.4546
$ws.Range("J46").Value = 889.6
$ws.Range("K46").Value = 1758.4546
$ws.Range("L46").Value = 889.6
$ws.Range("M46").Value = -1570.4546
$ws.Range("N46").Value = -1265.6
$ws.Range("H55").Value = 731.2778
$ws.Range("J55").Value = 1016.5714
$ws.Range("L55").Value = 1016.5714
$ws.Range("N55").Value = -1362.5714
$ws.Range("H95").Value = 29332.666
$ws.Range("J95").Value = 29332.666
$ws.Range("L95").Value = 29332.666
$ws.Range("N95").Value = -34824.666
$ws.Range("H140").Value = 40428.5
$ws.Range("J140").Value = 40428.5
$ws.Range("L140").Value = 40428.5
$ws.Range("N140").Value = -50788.5

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 11859.4
$ws.Range("I81").Value = 5489.7
$ws.Range("K81").Value = 10979.4
$ws.Range("M81").Value = -9918.4
$ws.Range("H84").Value = 11859.4
$ws.Range("I84").Value = 5489.7
$ws.Range("K84").Value = 54897
$ws.Range("M84").Value = -49593
$ws.Range("H97").Value = 22278.5
$ws.Range("J97").Value = 22278.5
$ws.Range("L97").Value = 22278.5
$ws.Range("N97").Value = -24260.5
$ws.Range("H122").Value = 5856.375
$ws.Range("I122").Value = 5621.5713
$ws.Range("K122").Value = 16864.7139
$ws.Range("M122").Value = -14414.7139
$ws.Range("H124").Value = 55808.223
$ws.Range("J124").Value = 55808.223
$ws.Range("L124").Value = 55808.223
$ws.Range("N124").Value = -65628.223
$ws.Range("H132").Value = 26040.936
$ws.Range("I132").Value = 24699
$ws.Range("K132").Value = 74097
$ws.Range("M132").Value = -71567
$ws.Range("H136").Value = 3407.9688
$ws.Range("I136").Value = 3744.5
$ws.Range("K136").Value = 11233.5
$ws.Range("M136").Value = -8683.5
